# Updated cryptos list on Sun Jun  9 07:26:34 UTC 2024 with GitHub Actions
#
# Applies the per-cell price/volume refresh (and the dogwifhat /
# InjectiveProtocol row swap) described by the diff against the
# "cryptos" worksheet.
#
# NumberFormat is forced to "@" (Text) before each write so that values
# such as "1.00", "0.0906", "69.383.07", etc. are stored verbatim as
# text (matching the workbook's existing inlineStr/Text cells) instead
# of being silently reinterpreted as numbers by Excel's smart-entry
# parser (which would otherwise drop things like the trailing zero in
# "1.00" or collapse "69.383.07" style values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "69.383.07"
Set-TextValue "E2" "  -0.11%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.688.82"
Set-TextValue "E3" "  -0.22%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.06%  "

# Row 5 - BNB
Set-TextValue "D5" "680.23"
Set-TextValue "E5" "  -2.01%  "

# Row 6 - Solana
Set-TextValue "D6" "159.44"
Set-TextValue "E6" "  -2.47%  "

# Row 7 - USDC
Set-TextValue "E7" "  +0.00%  "

# Row 8 - XRP
Set-TextValue "E8" "  -1.19%  "

# Row 9 - Dogecoin
Set-TextValue "E9" "  -1.61%  "

# Row 10 - Toncoin
Set-TextValue "D10" "7.06"
Set-TextValue "E10" "  -4.52%  "

# Row 11 - Cardano
Set-TextValue "E11" "  -1.80%  "

# Row 12 - ShibaInu
Set-TextValue "E12" "  -3.55%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "4.312.12"
Set-TextValue "E13" "  -0.05%  "

# Row 14 - Avalanche
Set-TextValue "D14" "32.48"
Set-TextValue "E14" "  -3.19%  "

# Row 15 - WrappedEther
Set-TextValue "D15" "3.691.04"
Set-TextValue "E15" "  -0.13%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "69.319.03"
Set-TextValue "E16" "  -0.24%  "

# Row 18 - Chainlink
Set-TextValue "E18" "  -1.72%  "

# Row 19 - Polkadot
Set-TextValue "D19" "6.43"
Set-TextValue "E19" "  -2.93%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "468.70"
Set-TextValue "E20" "  -3.06%  "

# Row 21 - Uniswap
Set-TextValue "D21" "10.01"
Set-TextValue "E21" "  +0.04%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.652"
Set-TextValue "E22" "  -2.38%  "

# Row 23 - Litecoin
Set-TextValue "D23" "79.92"
Set-TextValue "E23" "  -0.25%  "

# Row 24 - WrappedeETH
Set-TextValue "D24" "3.834.70"
Set-TextValue "E24" "  -0.04%  "

# Row 25 - Dai
Set-TextValue "E25" "  -0.01%  "

# Row 26 - PEPE
Set-TextValue "E26" "  -6.21%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "10.95"
Set-TextValue "E27" "  -4.44%  "

# Row 28 - RenderToken
Set-TextValue "D28" "9.15"
Set-TextValue "E28" "  -4.51%  "

# Row 29 - PancakeSwap
Set-TextValue "E29" "  -2.29%  "

# Row 30 - Fetch.AI
Set-TextValue "E30" "  -3.92%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "6.63"
Set-TextValue "E31" "  -3.81%  "

# Row 32 - ImmutableX
Set-TextValue "E32" "  -4.66%  "

# Row 33 - Binance-PegBSC-USD
Set-TextValue "D33" "1.00"
Set-TextValue "E33" "  -0.62%  "

# Row 34 - EthereumClassic
Set-TextValue "E34" "  -0.77%  "

# Row 35 - RenzoRestakedETH
Set-TextValue "D35" "3.677.37"
Set-TextValue "E35" "  +0.57%  "

# Row 36 - Kaspa
Set-TextValue "E36" "  -5.10%  "

# Row 37 - Aptos
Set-TextValue "E37" "  -3.14%  "

# Row 38 - Filecoin
Set-TextValue "E38" "  -1.97%  "

# Row 39 - USDe
Set-TextValue "E39" "  +0.00%  "

# Row 40 - Stacks
Set-TextValue "E40" "  -2.85%  "

# Row 41 - FirstDigitalUSD
Set-TextValue "E41" "  +0.02%  "

# Row 42 - Hedera
Set-TextValue "D42" "0.0906"
Set-TextValue "E42" "  -3.29%  "

# Row 43 - Monero
Set-TextValue "D43" "171.78"
Set-TextValue "E43" "  +4.81%  "

# Row 44 - Mantle
Set-TextValue "D44" "0.944"
Set-TextValue "E44" "  -1.10%  "

# Row 46 / 47 - dogwifhat and InjectiveProtocol swap ranking positions
# (InjectiveProtocol moves up to rank 46, dogwifhat moves down to rank 47)
Set-TextValue "B46" "InjectiveProtocol"
Set-TextValue "C46" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D46" "28.28"
Set-TextValue "E46" "  -6.18%  "

Set-TextValue "B47" "dogwifhat"
Set-TextValue "C47" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D47" "2.71"
Set-TextValue "E47" "  -4.43%  "

# Row 48 - FLOKI
Set-TextValue "E48" "  -3.61%  "

# Row 49 - ONDO
Set-TextValue "E49" "  -4.86%  "

# Row 50 - SuiNetwork
Set-TextValue "E50" "  -4.51%  "

# Row 51 - Cosmos
Set-TextValue "E51" "  -2.93%  "
